# Locations - SAMPLE.xlsx : replace the placeholder "Name" column values
# (Dakota, Gregs, Dunkins, RJ Logistics, Tess, Flex, McGann, Earls, Moes)
# with real landmark/building names, adjust column A's width for the new
# longer text, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "Bellagio Hotel"
$ws.Range("A3").Value  = "The Getty"
$ws.Range("A4").Value  = "Flatiron"
$ws.Range("A5").Value  = "KOIN Center"
$ws.Range("A6").Value  = "The Parthenon"
$ws.Range("A7").Value  = "Olympia Theater"
$ws.Range("A8").Value  = "Space Needle"
$ws.Range("A9").Value  = "Coors Field"
$ws.Range("A10").Value = "Anson Mills"

# Widen column A to fit the new (longer) names.
$ws.Columns.Item(1).ColumnWidth = 15

# Move the selected cell, as in the edited workbook.
$ws.Range("L7").Select()
